$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------------
# Sheet1: rebuild as an 8-column "Import Test Case" table (was a 2-col list)
# ---------------------------------------------------------------------------
$ws1.Cells.Clear()

$headers = @("Test Case Name","Description","Priority","Pre-condition","Post-condition","Steps","Expected Result","Type")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws1.Cells.Item(1, $i + 1).Value = $headers[$i]
}

$ws1.Cells.Item(2, 1).Value = "Sample A"
$ws1.Cells.Item(2, 2).Value = "This is a sample Test case"
$ws1.Cells.Item(2, 3).Value = "medium"
$ws1.Cells.Item(2, 4).Value = "None"
$ws1.Cells.Item(2, 5).Value = "None"
$ws1.Cells.Item(2, 6).Value = "1. Go to Google`n2. Display the object"
$ws1.Cells.Item(2, 7).Value = "1. Display google websites`n2. Display search result"
$ws1.Cells.Item(2, 8).Value = "Manual"

# Seed formatting (bordered, thin black border already present on original
# cells) by copy/paste-special so the thin border style is re-used exactly,
# then layer the extra look on top.
$ws2.Range("A2").Copy()
$ws1.Range("A1:H2").PasteSpecial(-4122)

$hdr = $ws1.Range("A1:H1")
$hdr.Font.Bold = $true
$hdr.Interior.Color = 15773696
$hdr.HorizontalAlignment = -4108

$body = $ws1.Range("A2:E2")
$bodyRest = $ws1.Range("H2")
$body.HorizontalAlignment = -4108
$bodyRest.HorizontalAlignment = -4108

$wrapCells = $ws1.Range("F2:G2")
$wrapCells.HorizontalAlignment = -4131
$wrapCells.WrapText = $true

$ws1.Rows.Item(2).RowHeight = 30

# Column widths (best effort; engine snaps to its internal width grid)
$ws1.Range("B1").ColumnWidth = 37.16666666
$ws1.Range("C1").ColumnWidth = 13.66666666
$ws1.Range("D1").ColumnWidth = 18.5
$ws1.Range("E1").ColumnWidth = 15.33333333
$ws1.Range("F1").ColumnWidth = 26
$ws1.Range("G1").ColumnWidth = 24.33333333
$ws1.Range("H1").ColumnWidth = 17

# ---------------------------------------------------------------------------
# Sheet2: same 4-column table, content unchanged (just re-asserted)
# ---------------------------------------------------------------------------
$ws2.Cells.Item(1, 1).Value = "Step"
$ws2.Cells.Item(1, 2).Value = "Definition"
$ws2.Cells.Item(1, 3).Value = "Expected Result"
$ws2.Cells.Item(1, 4).Value = "Type"

$ws2.Cells.Item(2, 1).Value = 1
$ws2.Cells.Item(2, 2).Value = "Go to google"
$ws2.Cells.Item(2, 3).Value = "Get it on"
$ws2.Cells.Item(2, 4).Value = "manual"

$ws2.Cells.Item(3, 1).Value = 2
$ws2.Cells.Item(3, 2).Value = "Go to google"
$ws2.Cells.Item(3, 3).Value = "Get it on"
$ws2.Cells.Item(3, 4).Value = "manual"

$ws2.Cells.Item(4, 1).Value = 3
$ws2.Cells.Item(4, 2).Value = "Go to google"
$ws2.Cells.Item(4, 3).Value = "Get it on"
$ws2.Cells.Item(4, 4).Value = "manual"

# ---------------------------------------------------------------------------
# Selection / active-sheet bookkeeping: Sheet1 becomes the selected tab.
# ---------------------------------------------------------------------------
$ws1.Activate()
$ws1.Range("C18").Select()
